# Weekly update: insert a new price-record row for "Papa" (Asterix, 1a
# (guarda)) at Feria Lagunitas de Puerto Montt, pushing the existing
# records (rows 654..737) down by one row (to 655..738).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 654; this shifts rows
# 654-737 down to 655-738 and extends the sheet dimension accordingly.
$ws.Rows.Item(654).Insert()

# Populate the newly inserted row 654 with the new record.
$ws.Cells.Item(654, 1).Value = 4
$ws.Cells.Item(654, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(654, 3).Value = "Los Lagos"
$ws.Cells.Item(654, 4).Value = 45154
$ws.Cells.Item(654, 5).Value = 10
$ws.Cells.Item(654, 6).Value = 100114001
$ws.Cells.Item(654, 7).Value = "Papa"
$ws.Cells.Item(654, 8).Value = "Asterix"
$ws.Cells.Item(654, 9).Value = "1a (guarda)"
$ws.Cells.Item(654, 10).Value = 300
$ws.Cells.Item(654, 11).Value = 20000
$ws.Cells.Item(654, 12).Value = 20000
$ws.Cells.Item(654, 13).Value = 20000
$ws.Cells.Item(654, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(654, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(654, 16).Value = 800
$ws.Cells.Item(654, 17).Value = 25
$ws.Cells.Item(654, 18).Value = "Hortaliza"
